$wb = $excel.ActiveWorkbook

# --- Sheet "Nos" (Nodes) ---
$ws1 = $wb.Worksheets.Item("Nos")
$ws1.Range("A3").Value = 0.1
$ws1.Range("B3").Value = 0
$ws1.Range("A4").Value = 0.1
$ws1.Range("B4").Value = 0.1
$ws1.Range("A5").Value = 0
$ws1.Range("B5").Value = 0.1

# --- Sheet "Incidencia" (Incidence) ---
$ws2 = $wb.Worksheets.Item("Incidencia")
$ws2.Range("B4").Value = 4

# New rows 5 and 6, formatted like row 3 (same style pattern as target)
$ws2.Range("A3:D3").Copy()
$ws2.Range("A5:D5").PasteSpecial(-4122)
$ws2.Range("A3:D3").Copy()
$ws2.Range("A6:D6").PasteSpecial(-4122)

$ws2.Range("A5").Value = 4
$ws2.Range("B5").Value = 1
$ws2.Range("C5").Value = 210000000000
$ws2.Range("D5").Value = 0.0002
$ws2.Range("A6").Value = 4
$ws2.Range("B6").Value = 2
$ws2.Range("C6").Value = 210000000000
$ws2.Range("D6").Value = 0.0002

# --- Sheet "Carregamento" (Loading) ---
$ws3 = $wb.Worksheets.Item("Carregamento")
$ws3.Range("A2").Value = 2
$ws3.Range("C2").Value = -100

# --- Sheet "Restricao" (Restriction) ---
$ws4 = $wb.Worksheets.Item("Restricao")

# Row 3: reformat to the "horizontal-center-only" style (like Incidencia!A2)
$ws2.Range("A2").Copy()
$ws4.Range("A3:B3").PasteSpecial(-4122)
$ws4.Range("A3").Value = 1
$ws4.Range("B3").Value = 2

# Row 4: values change, style (horizontal+vertical center) stays as-is
$ws4.Range("A4").Value = 4
$ws4.Range("B4").Value = 1

# Row 5: new values with the "horizontal-center-only" style
$ws2.Range("A2").Copy()
$ws4.Range("A5:B5").PasteSpecial(-4122)
$ws4.Range("A5").Value = 4
$ws4.Range("B5").Value = 2

# Rows 6 and 7: A/B cells removed entirely (clear content + formatting)
$ws4.Range("A6").Clear()
$ws4.Range("B6").Clear()
$ws4.Range("A7").Clear()
$ws4.Range("B7").Clear()

# --- Selections matching final cursor positions per sheet ---
$ws1.Range("D6").Select()
$ws2.Range("A6").Select()
$ws3.Range("C7").Select()
$ws4.Activate()
$ws4.Range("D8").Select()
$excel.ActiveWindow.Zoom = 96
